# Update "F" column (想去人数 / want-to-go count) values on the
# "展览" and "全部类型" worksheets to match the refreshed data export.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1)
$ws1.Range("F8").Value  = 4928
$ws1.Range("F9").Value  = 4928
$ws1.Range("F16").Value = 7984
$ws1.Range("F17").Value = 7984
$ws1.Range("F21").Value = 2465
$ws1.Range("F24").Value = 23
$ws1.Range("F25").Value = 2507
$ws1.Range("F34").Value = 6681
$ws1.Range("F40").Value = 104
$ws1.Range("F43").Value = 2493

# 全部类型 sheet (sheet4)
$ws4.Range("F11").Value = 4928
$ws4.Range("F12").Value = 4928
$ws4.Range("F18").Value = 7984
$ws4.Range("F19").Value = 7984
$ws4.Range("F23").Value = 2465
$ws4.Range("F28").Value = 2507
$ws4.Range("F38").Value = 6681
$ws4.Range("F41").Value = 104
$ws4.Range("F43").Value = 2493
